$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the mobile number in A2 (was 7708273167, now 7709856398)
$ws.Range("A2").Value = "7709856398"

# Restore the view: scroll to show A1 and select A2
$ws.Range("A2").Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
